# isa.study.xlsx / "Growth" sheet cleanup
#
# The "Term Source REF (DPBO:0000007)" and "Term Accession Number
# (DPBO:0000007)" columns (I2:J7) only ever held the placeholder text
# "user-specific" in every data row. Clear that stale placeholder data out
# of the table body - this also drops "user-specific" from the shared
# string table (it becomes unused) and tightens the sheet's used range back
# down to the real table extent (A1:K7) instead of the stale A1:DP7.
#
# Also update the saved selection on that sheet to C15 (previously J20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Growth")

$ws.Range("I2:J7").ClearContents()

$ws.Range("C15").Select()
